$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 227, shifting rows 227:234 down to 228:235
$ws.Rows.Item(227).Insert()

# Populate the new row 227 with the new "Pink Delight" entry.
# Common columns match the surrounding rows (A, B, C, E, F, G, H, I, J, R).
$ws.Cells.Item(227, 1).Value = 10
$ws.Cells.Item(227, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(227, 3).Value = "La Araucanía"
$ws.Cells.Item(227, 4).Value = 44706
$ws.Cells.Item(227, 4).NumberFormat = $ws.Cells.Item(228, 4).NumberFormat
$ws.Cells.Item(227, 5).Value = 9
$ws.Cells.Item(227, 6).Value = "Fruta"
$ws.Cells.Item(227, 7).Value = 100103
$ws.Cells.Item(227, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(227, 9).Value = 100103002
$ws.Cells.Item(227, 10).Value = "Ciruela"
$ws.Cells.Item(227, 11).Value = "Pink Delight"
$ws.Cells.Item(227, 12).Value = "Primera"
$ws.Cells.Item(227, 13).Value = 65
$ws.Cells.Item(227, 14).Value = 14000
$ws.Cells.Item(227, 15).Value = 14000
$ws.Cells.Item(227, 16).Value = 14000
$ws.Cells.Item(227, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(227, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(227, 19).Value = 778
$ws.Cells.Item(227, 20).Value = 18
